# Update sales figures on Sheet1 ("Sales July - Dec" column, D) for a
# handful of salespeople, per the source data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value  = 60000   # Andy Bernard
$ws.Range("D10").Value = 32000   # Phyllis Lapin
$ws.Range("D17").Value = 58500   # Hanna Moos
$ws.Range("D23").Value = 65700   # Patricio Simpson
$ws.Range("D36").Value = 51500   # Peter Franken
$ws.Range("D42").Value = 58300   # André Fonseca
$ws.Range("D43").Value = 62400   # Howard Snyder
